$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1270.4375
$ws.Range("I17").Value = 700
$ws.Range("J17").Value = 1308.4667
$ws.Range("K17").Value = 2100
$ws.Range("L17").Value = 3925.4001
$ws.Range("M17").Value = -1932
$ws.Range("N17").Value = -4261.4001
$ws.Range("H125").Value = 1433.8572
$ws.Range("I125").Value = 1236
$ws.Range("K125").Value = 11124
$ws.Range("M125").Value = -8664
$ws.Range("H137").Value = 3797.5
$ws.Range("I137").Value = 2924.8667
$ws.Range("J137").Value = 4804.385
$ws.Range("K137").Value = 8774.6001
$ws.Range("L137").Value = 14413.155
$ws.Range("M137").Value = -6224.6001
$ws.Range("N137").Value = -19513.155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7145766.5
$ws.Range("I32").Value = 7578154.5
$ws.Range("J32").Value = 11357.75
$ws.Range("K32").Value = 7578154.5
$ws.Range("L32").Value = 11357.75
$ws.Range("M32").Value = -7577867.5
$ws.Range("N32").Value = -11931.75
$ws.Range("H74").Value = 9625094
$ws.Range("I74").Value = 13891485
$ws.Range("J74").Value = 25712.375
$ws.Range("K74").Value = 13891485
$ws.Range("L74").Value = 25712.375
$ws.Range("M74").Value = -13890611
$ws.Range("N74").Value = -27460.375
$ws.Range("H77").Value = 9625094
$ws.Range("I77").Value = 13891485
$ws.Range("J77").Value = 25712.375
$ws.Range("K77").Value = 69457425
$ws.Range("L77").Value = 128561.875
$ws.Range("M77").Value = -69453057
$ws.Range("N77").Value = -137297.875
$ws.Range("H86").Value = 10314
$ws.Range("J86").Value = 10314
$ws.Range("L86").Value = 10314
$ws.Range("N86").Value = -12686
$ws.Range("H89").Value = 10314
$ws.Range("J89").Value = 10314
$ws.Range("L89").Value = 30942
$ws.Range("N89").Value = -42798
$ws.Range("H105").Value = 61786.668
$ws.Range("J105").Value = 77680
$ws.Range("L105").Value = 77680
$ws.Range("N105").Value = -84668
$ws.Range("H132").Value = 3399.5293
$ws.Range("I132").Value = 1316.2609
$ws.Range("K132").Value = 3948.7827
$ws.Range("M132").Value = -1418.7827

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 771.7
$ws.Range("I22").Value = 801.8889
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 801.8889
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -628.8889
$ws.Range("N22").Value = -846
$ws.Range("H80").Value = 1653.6666
$ws.Range("J80").Value = 483.125
$ws.Range("L80").Value = 483.125
$ws.Range("N80").Value = -2479.125
$ws.Range("H83").Value = 1653.6666
$ws.Range("J83").Value = 483.125
$ws.Range("L83").Value = 2415.625
$ws.Range("N83").Value = -12399.625
$ws.Range("H94").Value = 845.05554
$ws.Range("I94").Value = 877.25
$ws.Range("K94").Value = 877.25
$ws.Range("M94").Value = -426.25
$ws.Range("H99").Value = 7901
$ws.Range("I99").Value = 10383.546
$ws.Range("K99").Value = 10383.546
$ws.Range("M99").Value = -8885.546
$ws.Range("H105").Value = 2240.6924
$ws.Range("I105").Value = 1762.25
$ws.Range("K105").Value = 1762.25
$ws.Range("M105").Value = -15.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 671805
$ws.Range("I31").Value = 6050.6665
$ws.Range("J31").Value = 902258.4399999999
$ws.Range("K31").Value = 6050.6665
$ws.Range("L31").Value = 902258.4399999999
$ws.Range("M31").Value = -5755.6665
$ws.Range("N31").Value = -902848.4399999999
$ws.Range("H34").Value = 671805
$ws.Range("I34").Value = 6050.6665
$ws.Range("J34").Value = 902258.4399999999
$ws.Range("K34").Value = 6050.6665
$ws.Range("L34").Value = 902258.4399999999
$ws.Range("M34").Value = -5848.6665
$ws.Range("N34").Value = -902662.4399999999
$ws.Range("H132").Value = 4630.6665
$ws.Range("I132").Value = 3507.1428
$ws.Range("K132").Value = 10521.4284
$ws.Range("M132").Value = -7991.428400000001
$ws.Range("H134").Value = 3335633
$ws.Range("I134").Value = 5001949.5
$ws.Range("K134").Value = 15005848.5
$ws.Range("M134").Value = -15003313.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 4537.2856
$ws.Range("I33").Value = 3552.2
$ws.Range("K33").Value = 21313.2
$ws.Range("M33").Value = -21030.2
$ws.Range("H34").Value = 760
$ws.Range("I34").Value = 266.66666
$ws.Range("J34").Value = 1500
$ws.Range("K34").Value = 799.9999799999999
$ws.Range("L34").Value = 4500
$ws.Range("M34").Value = -715.9999799999999
$ws.Range("N34").Value = -4668
$ws.Range("H39").Value = 17270
$ws.Range("J39").Value = 17270
$ws.Range("L39").Value = 51810
$ws.Range("N39").Value = -52398
$ws.Range("H92").Value = 557299.5600000001
$ws.Range("I92").Value = 1112877.1
$ws.Range("J92").Value = 1722
$ws.Range("K92").Value = 3338631.3
$ws.Range("L92").Value = 5166
$ws.Range("M92").Value = -3337383.3
$ws.Range("N92").Value = -7662
$ws.Range("H132").Value = 1508.0741
$ws.Range("I132").Value = 1658.1538
$ws.Range("J132").Value = 1368.7142
$ws.Range("K132").Value = 14923.3842
$ws.Range("L132").Value = 12318.4278
$ws.Range("M132").Value = -12393.3842
$ws.Range("N132").Value = -17378.4278

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 87386780
$ws.Range("I11").Value = 147653120
$ws.Range("J11").Value = 7031668
$ws.Range("K11").Value = 147653120
$ws.Range("L11").Value = 7031668
$ws.Range("M11").Value = -147652981
$ws.Range("N11").Value = -7031946
$ws.Range("H39").Value = 74950
$ws.Range("J39").Value = 74950
$ws.Range("L39").Value = 74950
$ws.Range("N39").Value = -76014
$ws.Range("H69").Value = 111999
$ws.Range("J69").Value = 111999
$ws.Range("L69").Value = 111999
$ws.Range("N69").Value = -113497
$ws.Range("H72").Value = 111999
$ws.Range("J72").Value = 111999
$ws.Range("L72").Value = 335997
$ws.Range("N72").Value = -343485
$ws.Range("H74").Value = 44999
$ws.Range("J74").Value = 44999
$ws.Range("L74").Value = 44999
$ws.Range("N74").Value = -46871
$ws.Range("H77").Value = 44999
$ws.Range("J77").Value = 44999
$ws.Range("L77").Value = 134997
$ws.Range("N77").Value = -144357
$ws.Range("H95").Value = 166702050
$ws.Range("J95").Value = 166702050
$ws.Range("L95").Value = 166702050
$ws.Range("N95").Value = -166707542
$ws.Range("H132").Value = 90919576
$ws.Range("I132").Value = 100001530
$ws.Range("K132").Value = 300004590
$ws.Range("M132").Value = -300002060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1255
$ws.Range("I9").Value = 1255
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1255
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -1031
$ws.Range("N9").ClearContents()
$ws.Range("H13").Value = 3406.6667
$ws.Range("H20").Value = 32273.273
$ws.Range("J20").Value = 42501
$ws.Range("L20").Value = 42501
$ws.Range("N20").Value = -42953
$ws.Range("H24").Value = 70006.75
$ws.Range("I24").Value = 40006
$ws.Range("K24").Value = 40006
$ws.Range("M24").Value = -39663
$ws.Range("H101").Value = 67747.5
$ws.Range("J101").Value = 67747.5
$ws.Range("L101").Value = 67747.5
$ws.Range("N101").Value = -74237.5
$ws.Range("H106").Value = 41727
$ws.Range("J106").Value = 41727
$ws.Range("L106").Value = 41727
$ws.Range("N106").Value = -44251
$ws.Range("H122").Value = 6618.7095
$ws.Range("I122").Value = 4649.2856
$ws.Range("K122").Value = 13947.8568
$ws.Range("M122").Value = -11497.8568
$ws.Range("H136").Value = 45117.656
$ws.Range("I136").Value = 5972.6
$ws.Range("K136").Value = 17917.8
$ws.Range("M136").Value = -15367.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 267499.5
$ws.Range("I21").Value = 35000
$ws.Range("J21").Value = 499999
$ws.Range("K21").Value = 35000
$ws.Range("L21").Value = 499999
$ws.Range("M21").Value = -34765
$ws.Range("N21").Value = -500469
$ws.Range("H24").Value = 71257.5
$ws.Range("J24").Value = 71257.5
$ws.Range("L24").Value = 71257.5
$ws.Range("N24").Value = -71717.5
$ws.Range("H32").Value = 22249.5
$ws.Range("I32").Value = 22249.5
$ws.Range("K32").Value = 22249.5
$ws.Range("M32").Value = -21932.5
$ws.Range("H35").Value = 267499.5
$ws.Range("I35").Value = 35000
$ws.Range("J35").Value = 499999
$ws.Range("K35").Value = 35000
$ws.Range("L35").Value = 499999
$ws.Range("M35").Value = -34710
$ws.Range("N35").Value = -500579
$ws.Range("H100").Value = 2055.6
$ws.Range("I100").Value = 2061.7778
$ws.Range("K100").Value = 4123.5556
$ws.Range("M100").Value = -3582.5556
